$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fill in attendance hours ("B" column) for rows 74-85 that were left
# blank before - these represent missed / corrected daily entries.
$ws.Range("B74").Value = 2
$ws.Range("B75").Value = 0
$ws.Range("B76").Value = 0
$ws.Range("B77").Value = 1
$ws.Range("B78").Value = 0
$ws.Range("B79").Value = 0
$ws.Range("B80").Value = 3
$ws.Range("B81").Value = 0
$ws.Range("B82").Value = 12
$ws.Range("B83").Value = 10
$ws.Range("B84").Value = 13
$ws.Range("B85").Value = 14

# Move the active selection to reflect where the editor left off.
$ws.Range("H80").Select()
